# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# updates to the Louisoix_Profits workbook, as described by the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 2420.75
$ws.Range("J19").Value = 1671
$ws.Range("L19").Value = 1671
$ws.Range("N19").Value = -2021
# Row 28
$ws.Range("H28").Value = 273.26666
$ws.Range("I28").Value = 265.63635
$ws.Range("K28").Value = 265.63635
$ws.Range("M28").Value = 219.36365
# Row 32
$ws.Range("H32").Value = 2380.6428
$ws.Range("I32").Value = 1426.6666
$ws.Range("J32").Value = 2640.818
$ws.Range("K32").Value = 1426.6666
$ws.Range("L32").Value = 2640.818
$ws.Range("M32").Value = -1100.6666
$ws.Range("N32").Value = -3292.818
# Row 53
$ws.Range("H53").Value = 465.81818
$ws.Range("I53").Value = 279.375
$ws.Range("J53").Value = 963
$ws.Range("K53").Value = 279.375
$ws.Range("L53").Value = 963
$ws.Range("M53").Value = 357.625
$ws.Range("N53").Value = -2237
# Row 70
$ws.Range("H70").Value = 1766.5555
$ws.Range("I70").Value = 1414.2858
$ws.Range("J70").Value = 2999.5
$ws.Range("K70").Value = 4242.857400000001
$ws.Range("L70").Value = 8998.5
$ws.Range("M70").Value = -3972.857400000001
$ws.Range("N70").Value = -9538.5
# Row 73
$ws.Range("H73").Value = 1766.5555
$ws.Range("I73").Value = 1414.2858
$ws.Range("J73").Value = 2999.5
$ws.Range("K73").Value = 4242.857400000001
$ws.Range("L73").Value = 8998.5
$ws.Range("M73").Value = -3306.857400000001
$ws.Range("N73").Value = -10870.5
# Row 80
$ws.Range("H80").Value = 26642.783
$ws.Range("I80").Value = 54852.273
$ws.Range("K80").Value = 164556.819
$ws.Range("M80").Value = -163558.819
# Row 83
$ws.Range("H83").Value = 26642.783
$ws.Range("I83").Value = 54852.273
$ws.Range("K83").Value = 493670.457
$ws.Range("M83").Value = -488678.457
# Row 98
$ws.Range("H98").Value = 1132
$ws.Range("I98").Value = 1143.4546
$ws.Range("K98").Value = 1143.4546
$ws.Range("M98").Value = 354.5454
# Row 112
$ws.Range("H112").Value = 2288.9375
$ws.Range("J112").Value = 2305.7693
$ws.Range("L112").Value = 6917.3079
$ws.Range("N112").Value = -9133.3079
# Row 116
$ws.Range("H116").Value = 5738.722
$ws.Range("J116").Value = 5933.125
$ws.Range("L116").Value = 5933.125
$ws.Range("N116").Value = -12817.125
# Row 122
$ws.Range("H122").Value = 1132
$ws.Range("I122").Value = 1143.4546
$ws.Range("K122").Value = 3430.3638
$ws.Range("M122").Value = -980.3638000000001
# Row 127
$ws.Range("H127").Value = 1572.0834
$ws.Range("I127").Value = 1535.5555
$ws.Range("K127").Value = 4606.666499999999
$ws.Range("M127").Value = 353.3335000000006
# Row 129
$ws.Range("H129").Value = 607.25
$ws.Range("I129").Value = 607.25
$ws.Range("K129").Value = 1821.75
$ws.Range("M129").Value = 3178.25
# Row 131
$ws.Range("H131").Value = 23257.666
$ws.Range("I131").Value = 1866.4
$ws.Range("K131").Value = 5599.200000000001
$ws.Range("M131").Value = -559.2000000000007
# Row 132
$ws.Range("H132").Value = 2999.9614
$ws.Range("I132").Value = 2999.9614
$ws.Range("K132").Value = 8999.8842
$ws.Range("M132").Value = -6469.8842
# Row 137
$ws.Range("H137").Value = 2396.7334
$ws.Range("I137").Value = 1109.1666
$ws.Range("K137").Value = 3327.4998
$ws.Range("M137").Value = -777.4998000000001
# Row 141
$ws.Range("H141").Value = 7849.6665
$ws.Range("I141").Value = 9119.6
$ws.Range("K141").Value = 27358.8
$ws.Range("M141").Value = -22178.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4606.067
$ws.Range("I45").Value = 3485.375
$ws.Range("K45").Value = 3485.375
$ws.Range("M45").Value = -3108.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3130.4285
$ws.Range("I134").Value = 2139.8
$ws.Range("K134").Value = 6419.400000000001
$ws.Range("M134").Value = -3884.400000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 39290.184
$ws.Range("I86").Value = 47271.332
$ws.Range("K86").Value = 47271.332
$ws.Range("M86").Value = -46148.332
# Row 89
$ws.Range("H89").Value = 39290.184
$ws.Range("I89").Value = 47271.332
$ws.Range("K89").Value = 236356.66
$ws.Range("M89").Value = -230740.66
# Row 122
$ws.Range("H122").Value = 2179.6667
$ws.Range("I122").Value = 1995.6
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 5986.799999999999
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -3536.799999999999
$ws.Range("N122").Value = -14200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 10374.75
$ws.Range("I3").Value = 7166.3335
$ws.Range("K3").Value = 21499.0005
$ws.Range("M3").Value = -21387.0005
# Row 40
$ws.Range("H40").Value = 137.36363
$ws.Range("I40").Value = 111
$ws.Range("K40").Value = 444
$ws.Range("M40").Value = -375
# Row 129
$ws.Range("H129").Value = 718541.2
$ws.Range("I129").Value = 11859.5
$ws.Range("K129").Value = 35578.5
$ws.Range("M129").Value = -30578.5
# Row 131
$ws.Range("H131").Value = 3852382.5
$ws.Range("I131").Value = 27550.75
$ws.Range("J131").Value = 4547806.5
$ws.Range("K131").Value = 82652.25
$ws.Range("L131").Value = 13643419.5
$ws.Range("M131").Value = -77612.25
$ws.Range("N131").Value = -13653499.5
# Row 137
$ws.Range("H137").Value = 4690
$ws.Range("I137").Value = 2397.5
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 7192.5
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -2092.5
$ws.Range("N137").Value = -28200
# Row 139
$ws.Range("H139").Value = 1518.375
$ws.Range("I139").Value = 1518.375
$ws.Range("K139").Value = 4555.125
$ws.Range("M139").Value = 584.875
# Row 140
$ws.Range("H140").Value = 2118.8696
$ws.Range("I140").Value = 1720.7142
$ws.Range("K140").Value = 5162.142599999999
$ws.Range("M140").Value = 17.85740000000078

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 70
$ws.Range("H70").Value = 13229.083
$ws.Range("I70").Value = 12861.111
$ws.Range("J70").Value = 14333
$ws.Range("K70").Value = 12861.111
$ws.Range("L70").Value = 14333
$ws.Range("M70").Value = -12591.111
$ws.Range("N70").Value = -14873
# Row 73
$ws.Range("H73").Value = 13229.083
$ws.Range("I73").Value = 12861.111
$ws.Range("J73").Value = 14333
$ws.Range("K73").Value = 12861.111
$ws.Range("L73").Value = 14333
$ws.Range("M73").Value = -11925.111
$ws.Range("N73").Value = -16205
# Row 132
$ws.Range("H132").Value = 33064.125
$ws.Range("I132").Value = 38556.258
$ws.Range("K132").Value = 115668.774
$ws.Range("M132").Value = -113138.774

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 36519.875
$ws.Range("I22").Value = 62125.61
$ws.Range("J22").Value = 3598.2144
$ws.Range("K22").Value = 62125.61
$ws.Range("L22").Value = 3598.2144
$ws.Range("M22").Value = -61830.61
$ws.Range("N22").Value = -4188.2144
# Row 27
$ws.Range("H27").Value = 36519.875
$ws.Range("I27").Value = 62125.61
$ws.Range("J27").Value = 3598.2144
$ws.Range("K27").Value = 62125.61
$ws.Range("L27").Value = 3598.2144
$ws.Range("M27").Value = -62018.61
$ws.Range("N27").Value = -3812.2144
# Row 46
$ws.Range("H46").Value = 13811.241
$ws.Range("I46").Value = 18999.111
$ws.Range("J46").Value = 5322
$ws.Range("K46").Value = 18999.111
$ws.Range("L46").Value = 5322
$ws.Range("M46").Value = -18811.111
$ws.Range("N46").Value = -5698
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
# Row 93
$ws.Range("H93").Value = 2029.5
$ws.Range("I93").Value = 1868.5
$ws.Range("K93").Value = 1868.5
$ws.Range("M93").Value = -620.5
# Row 132
$ws.Range("H132").Value = 31007.93
$ws.Range("I132").Value = 35828.918
$ws.Range("K132").Value = 107486.754
$ws.Range("M132").Value = -104956.754

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 100000
$ws.Range("I8").Value = 100000
$ws.Range("K8").Value = 100000
$ws.Range("M8").Value = -99860
# Row 96
$ws.Range("H96").Value = 3150
$ws.Range("I96").Value = 3150
$ws.Range("K96").Value = 3150
$ws.Range("M96").Value = -1777
# Row 108
$ws.Range("H108").Value = 59985
$ws.Range("J108").Value = 59985
$ws.Range("L108").Value = 59985
$ws.Range("N108").Value = -67665
